# Apply gene-name venn data updates and header label updates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1: replace numeric index headers with descriptive combination labels
$ws.Range("A1").Value = "DESeq2/LimmaVoom/EdgeR"
$ws.Range("B1").Value = "DESeq2/EdgeR"
$ws.Range("C1").Value = "DESeq2/LimmaVoom"
$ws.Range("D1").Value = "EdgeR/LimmaVoom"
$ws.Range("E1").Value = "DESeq2"
$ws.Range("F1").Value = "LimmaVoom"
$ws.Range("G1").Value = "EdgeR"

# Swap B2 / B3 (RPS9 <-> CROCCP2)
$ws.Range("B2").Value = "CROCCP2"
$ws.Range("B3").Value = "RPS9"

# Swap E23 / E24 (ZNF462 <-> RP11-420C9.1)
$ws.Range("E23").Value = "RP11-420C9.1"
$ws.Range("E24").Value = "ZNF462"

# Shift G27 -> G28 -> G29 -> G27 (NME2P1, PDE3A, CH507-528H12.1)
$ws.Range("G27").Value = "CH507-528H12.1"
$ws.Range("G28").Value = "NME2P1"
$ws.Range("G29").Value = "PDE3A"

# Swap F34 / F35 (SOCS1 <-> AC016757.3)
$ws.Range("F34").Value = "AC016757.3"
$ws.Range("F35").Value = "SOCS1"

# Swap F39 / F40 (NEO1 <-> CTD-2007H13.3)
$ws.Range("F39").Value = "CTD-2007H13.3"
$ws.Range("F40").Value = "NEO1"

# Shift F43 -> F44 -> F45 -> F43 (Y_RNA, ZNF252P-AS1, CTD-2186M15.3)
$ws.Range("F43").Value = "CTD-2186M15.3"
$ws.Range("F44").Value = "Y_RNA"
$ws.Range("F45").Value = "ZNF252P-AS1"

# Swap G53 / G54 (RP11-545I5.3 <-> RP11-531H8.2)
$ws.Range("G53").Value = "RP11-531H8.2"
$ws.Range("G54").Value = "RP11-545I5.3"
